$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 2278.5862
$ws.Range("I40").Value = 2160.389
$ws.Range("J40").Value = 2472
$ws.Range("K40").Value = 2160.389
$ws.Range("L40").Value = 2472
$ws.Range("M40").Value = -1985.389
$ws.Range("N40").Value = -2822

# Row 70: Consecrating Congregation
$ws.Range("H70").Value = 51841.19
$ws.Range("I70").Value = 3452.8572
$ws.Range("J70").Value = 76035.36
$ws.Range("K70").Value = 10358.5716
$ws.Range("L70").Value = 228106.08
$ws.Range("M70").Value = -10088.5716
$ws.Range("N70").Value = -228646.08

# Row 73: Curbing the Contagion (L)
$ws.Range("H73").Value = 51841.19
$ws.Range("I73").Value = 3452.8572
$ws.Range("J73").Value = 76035.36
$ws.Range("K73").Value = 10358.5716
$ws.Range("L73").Value = 228106.08
$ws.Range("M73").Value = -9422.571599999999
$ws.Range("N73").Value = -229978.08

# Row 101: Edge of the Arcane
$ws.Range("H101").Value = 25000624
$ws.Range("I101").Value = 33333916
$ws.Range("J101").Value = 750
$ws.Range("K101").Value = 100001748
$ws.Range("L101").Value = 2250
$ws.Range("M101").Value = -100000126
$ws.Range("N101").Value = -5494

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1148.7778
$ws.Range("I132").Value = 1148.7778
$ws.Range("K132").Value = 3446.3334
$ws.Range("M132").Value = -916.3334000000004

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1540.6666
$ws.Range("I137").Value = 1540.6666
$ws.Range("K137").Value = 4621.9998
$ws.Range("M137").Value = -2071.9998

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6942.067
$ws.Range("I32").Value = 7080.7856
$ws.Range("K32").Value = 7080.7856
$ws.Range("M32").Value = -6793.7856

# Row 41: Skillet Scandal
$ws.Range("H41").Value = 1995.6666
$ws.Range("I41").Value = 1995.6666
$ws.Range("K41").Value = 1995.6666
$ws.Range("M41").Value = -1581.6666

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 7213.4287
$ws.Range("J61").Value = 7213.4287
$ws.Range("L61").Value = 7213.4287
$ws.Range("N61").Value = -7637.4287

# Row 97: Ore for Me
$ws.Range("H97").Value = 1491.3334
$ws.Range("I97").Value = 1639.6
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 1639.6
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = -1143.6
$ws.Range("N97").Value = -1742

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 7213.4287
$ws.Range("J136").Value = 7213.4287
$ws.Range("L136").Value = 21640.2861
$ws.Range("N136").Value = -26740.2861

$ws = $wb.Worksheets.Item("CRP")
# Row 69: Landing the Big One
$ws.Range("H69").Value = 34444
$ws.Range("I69").Value = 34444
$ws.Range("K69").Value = 34444
$ws.Range("M69").Value = -33695

# Row 72: Fishing for Profits (L)
$ws.Range("H72").Value = 34444
$ws.Range("I72").Value = 34444
$ws.Range("K72").Value = 103332
$ws.Range("M72").Value = -99588

# Row 80: The Long Armillae of the Law
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -32246

# Row 83: Wooden Ambitions (L)
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -101232

# Row 99: O Pine
$ws.Range("H99").Value = 4989
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4989
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = $null
$ws.Range("M99").Value = 4989
$ws.Range("N99").Value = -7985

# Row 126: A Better Conductor
$ws.Range("H126").Value = 4989
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4989
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = $null
$ws.Range("M126").Value = 14967
$ws.Range("N126").Value = -19907

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2099.4443
$ws.Range("I132").Value = 2087.5
$ws.Range("J132").Value = 2195
$ws.Range("K132").Value = 6262.5
$ws.Range("L132").Value = 6585
$ws.Range("M132").Value = -3732.5
$ws.Range("N132").Value = -11645

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 105385.86
$ws.Range("I2").Value = 110077
$ws.Range("J2").Value = 101121.18
$ws.Range("K2").Value = 660462
$ws.Range("L2").Value = 606727.08
$ws.Range("M2").Value = -660349
$ws.Range("N2").Value = -606953.08

# Row 22: A Total Nut Job
$ws.Range("H22").Value = 2257.8333
$ws.Range("J22").Value = 2386.75
$ws.Range("L22").Value = 7160.25
$ws.Range("N22").Value = -7498.25

# Row 27: Brain Food
$ws.Range("H27").Value = 2257.8333
$ws.Range("J27").Value = 2386.75
$ws.Range("L27").Value = 7160.25
$ws.Range("N27").Value = -7364.25

# Row 46: Feeding Frenzy
$ws.Range("H46").Value = 674.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 674.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = $null
$ws.Range("M46").Value = 2024.25
$ws.Range("N46").Value = -2206.25

$ws = $wb.Worksheets.Item("ALC")
# Row 120: Supreme Official Strategy Guide
$ws.Range("H120").Value = 2000
$ws.Range("I120").Value = 2000
$ws.Range("K120").Value = 6000
$ws.Range("M120").Value = -1162

$ws = $wb.Worksheets.Item("CUL")
# Row 122: Salt of the North
$ws.Range("H122").Value = 896.3333
$ws.Range("J122").Value = 896.3333
$ws.Range("L122").Value = 8066.9997
$ws.Range("N122").Value = -12966.9997

# Row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 5317.316
$ws.Range("I140").Value = 1859.4286
$ws.Range("J140").Value = 14999.4
$ws.Range("K140").Value = 5578.2858
$ws.Range("L140").Value = 44998.2
$ws.Range("M140").Value = -398.2857999999997
$ws.Range("N140").Value = -55358.2

$ws = $wb.Worksheets.Item("GSM")
# Row 103: Ring in the New
$ws.Range("H103").Value = 46663.332
$ws.Range("J103").Value = 46663.332
$ws.Range("L103").Value = 46663.332
$ws.Range("N103").Value = -49007.332

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 7555
$ws.Range("J7").Value = 8768.462
$ws.Range("L7").Value = 8768.462
$ws.Range("N7").Value = -8992.462

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 5097.5
$ws.Range("I22").Value = 4847.125
$ws.Range("J22").Value = 5431.3335
$ws.Range("K22").Value = 4847.125
$ws.Range("L22").Value = 5431.3335
$ws.Range("M22").Value = -4552.125
$ws.Range("N22").Value = -6021.3335

# Row 27: Fire and Hide
$ws.Range("H27").Value = 5097.5
$ws.Range("I27").Value = 4847.125
$ws.Range("J27").Value = 5431.3335
$ws.Range("K27").Value = 4847.125
$ws.Range("L27").Value = 5431.3335
$ws.Range("M27").Value = -4740.125
$ws.Range("N27").Value = -5645.3335

# Row 43: Subordinate Clause
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = $null

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 3214.5
$ws.Range("I46").Value = 3595.6667
$ws.Range("J46").Value = 2833.3333
$ws.Range("K46").Value = 3595.6667
$ws.Range("L46").Value = 2833.3333
$ws.Range("M46").Value = -3407.6667
$ws.Range("N46").Value = -3209.3333

# Row 56: Hold On Tight
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = $null

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 2540.818
$ws.Range("I61").Value = 1935.8
$ws.Range("J61").Value = 3045
$ws.Range("K61").Value = 1935.8
$ws.Range("L61").Value = 3045
$ws.Range("M61").Value = -1733.8
$ws.Range("N61").Value = -3449

# Row 63: From Mud to Mourning
$ws.Range("H63").Value = 89538
$ws.Range("I63").Value = 89077
$ws.Range("K63").Value = 89077
$ws.Range("M63").Value = -88328

# Row 66: These Boots Are Made for Hawkin' (L)
$ws.Range("H66").Value = 89538
$ws.Range("I66").Value = 89077
$ws.Range("K66").Value = 267231
$ws.Range("M66").Value = -263487

# Row 74: Overall, We Blend In
$ws.Range("H74").Value = 58598.5
$ws.Range("I74").Value = 58598.5
$ws.Range("K74").Value = 58598.5
$ws.Range("M74").Value = -57600.5

# Row 77: Eviction Notice (L)
$ws.Range("H77").Value = 58598.5
$ws.Range("I77").Value = 58598.5
$ws.Range("K77").Value = 175795.5
$ws.Range("M77").Value = -170803.5

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 2679
$ws.Range("I82").Value = 3124
$ws.Range("J82").Value = 899
$ws.Range("K82").Value = 3124
$ws.Range("L82").Value = 899
$ws.Range("M82").Value = -2763
$ws.Range("N82").Value = -1621

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 2679
$ws.Range("I85").Value = 3124
$ws.Range("J85").Value = 899
$ws.Range("K85").Value = 3124
$ws.Range("L85").Value = 899
$ws.Range("M85").Value = -1876
$ws.Range("N85").Value = -3395

# Row 113: Peace in Rest
$ws.Range("H113").Value = 2540.818
$ws.Range("I113").Value = 1935.8
$ws.Range("J113").Value = 3045
$ws.Range("K113").Value = 1935.8
$ws.Range("L113").Value = 3045
$ws.Range("M113").Value = 234.2
$ws.Range("N113").Value = -7385

# Row 122: Hell on Leather
$ws.Range("H122").Value = 7837.552
$ws.Range("I122").Value = 8266.467000000001
$ws.Range("K122").Value = 24799.401
$ws.Range("M122").Value = -22349.401

# Row 126: Battered Books
$ws.Range("H126").Value = 7555
$ws.Range("J126").Value = 8768.462
$ws.Range("L126").Value = 26305.386
$ws.Range("N126").Value = -31245.386

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 2639.3076
$ws.Range("I132").Value = 2130.6
$ws.Range("K132").Value = 6391.799999999999
$ws.Range("M132").Value = -3861.799999999999

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4500.25

$ws = $wb.Worksheets.Item("WVR")
# Row 70: An Account of My Boots
$ws.Range("H70").Value = 32525
$ws.Range("I70").Value = 29995
$ws.Range("K70").Value = 29995
$ws.Range("M70").Value = -29680

# Row 73: Soot in My Hair and Scars on My Feet (L)
$ws.Range("H73").Value = 32525
$ws.Range("I73").Value = 29995
$ws.Range("K73").Value = 29995
$ws.Range("M73").Value = -28903

# Row 75: Storm upon Bald Mountain
$ws.Range("H75").Value = 91556.336
$ws.Range("I75").Value = 89559
$ws.Range("K75").Value = 89559
$ws.Range("M75").Value = -88623

# Row 78: Abrupt Apprentices (L)
$ws.Range("H78").Value = 91556.336
$ws.Range("I78").Value = 89559
$ws.Range("K78").Value = 268677
$ws.Range("M78").Value = -263997

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 948.8
$ws.Range("J122").Value = 897.5
$ws.Range("L122").Value = 2692.5
$ws.Range("N122").Value = -7592.5

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 4069.2632
$ws.Range("I126").Value = 2743.0833
$ws.Range("J126").Value = 6342.7144
$ws.Range("K126").Value = 8229.249899999999
$ws.Range("L126").Value = 19028.1432
$ws.Range("M126").Value = -5759.249899999999
$ws.Range("N126").Value = -23968.1432
